$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 236, shifting existing rows 236:293 down to 237:294
$ws.Rows("236:236").Insert()

# Populate the new row 236 with its values
$ws.Range("A236").Value = 11
$ws.Range("B236").Value = "Vega Monumental Concepción"
$ws.Range("C236").Value = "Bíobío"
$ws.Range("D236").Value = 44876
$ws.Range("D236").NumberFormat = $ws.Range("D237").NumberFormat
$ws.Range("E236").Value = 8
$ws.Range("F236").Value = 100114013
$ws.Range("G236").Value = "Zanahoria"
$ws.Range("H236").Value = "Sin especificar"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 270
$ws.Range("K236").Value = 15000
$ws.Range("L236").Value = 16000
$ws.Range("M236").Value = 15556
$ws.Range("N236").Value = '$/saco 20 kilos'
$ws.Range("O236").Value = "Región Metropolitana"
$ws.Range("P236").Value = 778
$ws.Range("Q236").Value = 20
$ws.Range("R236").Value = "Hortaliza"
